$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @"
questions = [
    {
        "title": "By default, the Ionic grid takes 100% width.Which of the following code can allow you to set a specific width based on the screen size?",
        "ques_type": 2,
        "options": [
            "&ltion-grid class=\"ion-grid-width\"&gt",
            "&ltion-grid size=\"10\" offset=\"2\"&gt",
            "&ltion-grid style=\"width: 540px\"&gt",
            "&ltion-grid fixed=\"true\"&gt"
        ],
        "score": "&ltion-grid fixed=\"true\"&gt"
    },
    {
        "title": "To publish your app as PWA, you wrote the following script in index.html. Finally, you copied [project_folder]/platforms/browser/www contents to your http server.What CLI command should you execute before copying them to the web server?&lt!--script&gt\n    if ('serviceWorker' in navigator) {\n      navigator.serviceWorker.register('service-worker.js')\n        .then(() =&gt console.log('service worker is installed!'))\n        .catch(err =&gt console.log('Error found during service worker installation', err))\n    }\n&lt/script--&gt",
        "ques_type": 2,
        "options": [
            "ionic cordova platform add browsernpm run ionic:build --prod",
            "ionic cordova platform add browserionic build browser --prod --release",
            "ionic cordova build androidionic cordova emulate android",
            "ionic cordova build pwaionic cordova emulate pwa"
        ],
        "score": "ionic cordova platform add browserionic build browser --prod --release"
    },
    {
        "title": "In Ionic Framework 5+, you declare the providers array in app.module.ts as shown in the code below.How can you write code in app.module.ts to import SQLite?providers: [\n StatusBar,\n SplashScreen,\n SQLite,\n { provide: RouteReuseStrategy, useClass: IonicRouteStrategy }",
        "ques_type": 2,
        "options": [
            "import { SQLite, SQLiteObject } from '@ionic-native/sqlite'",
            "import { SQLite } from '@ionic-native/sqlite'",
            "import { SQLite } from '@ionic-native/sqlite/ngx'",
            "import { SQLite, SQLiteObject } from '@ionic-native/sqlite/ngx'"
        ],
        "score": "import { SQLite } from '@ionic-native/sqlite/ngx'"
    },
    {
        "title": "When you open src/index.tsx in a React project, you can see the code shown below.Which of the following statements are true about this code?import React from 'react'\nimport ReactDOM from 'react-dom'\nimport App from './App'\nReactDOM.render(&ltApp /&gt, document.getElementById('root'))",
        "ques_type": 15,
        "options": [
            "The first line is not pulling in some dependencies.",
            "The first line allows you to write components in an HTML-like syntax called JSX.",
            "The ReactDOM.render method is not the browser-specific way to take your component.",
            "The ReactDOM.render method is the way to render components to a specific DOM node.",
            "The third line imports an optional component for your app."
        ],
        "score": [
            "The first line allows you to write components in an HTML-like syntax called JSX.",
            "The ReactDOM.render method is the way to render components to a specific DOM node."
        ]
    }
]
"@

$ws.Range("A2").EntireRow.Delete()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows.Item(1).AutoFit()
